# edit.ps1 - apply the Projeto_Integrador_V2.docx changes described by the diff
#
# Summary of changes:
#  1. Paragraph "JSON" -> "Conexão entre o app e o BD: ConnectionFactory.java"
#     (two runs) and gains the _GoBack bookmark (moved from further down).
#  2. Paragraph "Modelo Mobile(Renata)" -> "Template de app(Renata)" with
#     proofErr spell/grammar markers.
#  3. A brand-new paragraph "CRUDs(Rodrigo)" (same list style) is inserted
#     right after it, also with proofErr markers.
#  4. The old _GoBack bookmark that used to sit after "...dentro da equipe"
#     is removed (it moved to change #1).
#  5. <w:lastRenderedPageBreak/> is added to the run that starts "O
#     planejamento de teste..." and removed from the "Definição:" run that
#     immediately follows it.

$d = $word.ActiveDocument

function Find-ParaIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 1+4) "O planejamento de teste..." gains a lastRenderedPageBreak, and the
#      following "Definição:" paragraph loses its lastRenderedPageBreak.
# ---------------------------------------------------------------------
$iPlan = Find-ParaIndex $d "O planejamento de teste dever seguir o modelo que vai ser disponível no AVA da disciplina Teste/Qualidade de Software."
$pPlan = $d.Paragraphs.Item($iPlan)
$pDef  = $d.Paragraphs.Item($iPlan + 1)
$rngPlanDef = $d.Range($pPlan.Range.Start, $pDef.Range.End)

$xmlPlanDef = $pkgHeader +
    '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>O planejamento de teste dever seguir o modelo que vai ser disponível no AVA da disciplina Teste/Qualidade de Software.</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Definição:</w:t></w:r></w:p>' +
    $pkgFooter
[void]$rngPlanDef.InsertXML($xmlPlanDef)

# ---------------------------------------------------------------------
# 2) Remove the old _GoBack bookmark after "...dentro da equipe".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3) "Modelo Mobile(Renata)" -> "Template de app(Renata)" (+ proofErr) and
#    insert new "CRUDs(Rodrigo)" paragraph right after it.
# ---------------------------------------------------------------------
$iModelo = Find-ParaIndex $d "Modelo Mobile(Renata)"
$pModelo = $d.Paragraphs.Item($iModelo)

$xmlModeloCruds = $pkgHeader +
    '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Template</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>app</w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Renata)</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>CRUDs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Rodrigo)</w:t></w:r></w:p>' +
    $pkgFooter
[void]$pModelo.Range.InsertXML($xmlModeloCruds)

# ---------------------------------------------------------------------
# 4) "JSON" -> "Conexão entre o app e o BD: ConnectionFactory.java" plus the
#    _GoBack bookmark moved here (as the last content in the paragraph).
# ---------------------------------------------------------------------
$iJson = Find-ParaIndex $d "JSON"
$pJson = $d.Paragraphs.Item($iJson)

$xmlJson = $pkgHeader +
    '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>Conexão entre o app e o BD:</w:t></w:r><w:r><w:t xml:space="preserve"> ConnectionFactory.java</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
    $pkgFooter
[void]$pJson.Range.InsertXML($xmlJson)

Write-Output "done"
